$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 83.63636
$ws.Range("I2").Value = 83.63636
$ws.Range("K2").Value = 83.63636
$ws.Range("M2").Value = 29.36364

$ws.Range("H62").Value = 6948667.5
$ws.Range("I62").Value = 11366549
$ws.Range("J62").Value = 6282.4287
$ws.Range("K62").Value = 11366549
$ws.Range("L62").Value = 6282.4287
$ws.Range("M62").Value = -11365925
$ws.Range("N62").Value = -7530.4287

$ws.Range("H65").Value = 6948667.5
$ws.Range("I65").Value = 11366549
$ws.Range("J65").Value = 6282.4287
$ws.Range("K65").Value = 56832745
$ws.Range("L65").Value = 31412.1435
$ws.Range("M65").Value = -56829625
$ws.Range("N65").Value = -37652.14350000001

$ws.Range("H106").Value = 2682.5557
$ws.Range("I106").Value = 2546.5293
$ws.Range("K106").Value = 2546.5293
$ws.Range("M106").Value = -1915.5293

$ws.Range("H107").Value = 64473.438
$ws.Range("I107").Value = 102028.7
$ws.Range("J107").Value = 1881.3334
$ws.Range("K107").Value = 102028.7
$ws.Range("L107").Value = 1881.3334
$ws.Range("M107").Value = -100108.7
$ws.Range("N107").Value = -5721.3334

$ws.Range("H124").Value = 60769.6
$ws.Range("J124").Value = 60769.6
$ws.Range("L124").Value = 60769.6
$ws.Range("N124").Value = -70589.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1316.8334
$ws.Range("I2").Value = 1377.6
$ws.Range("J2").Value = 1013
$ws.Range("K2").Value = 1377.6
$ws.Range("L2").Value = 1013
$ws.Range("M2").Value = -1264.6
$ws.Range("N2").Value = -1239

$ws.Range("H32").Value = 3890.092
$ws.Range("I32").Value = 3890.092
$ws.Range("K32").Value = 3890.092
$ws.Range("M32").Value = -3603.092

$ws.Range("H61").Value = 2401.7778
$ws.Range("I61").Value = 2396
$ws.Range("K61").Value = 2396
$ws.Range("M61").Value = -2184

$ws.Range("H74").Value = 1424.6052
$ws.Range("I74").Value = 1354.5
$ws.Range("K74").Value = 1354.5
$ws.Range("M74").Value = -480.5

$ws.Range("H77").Value = 1424.6052
$ws.Range("I77").Value = 1354.5
$ws.Range("K77").Value = 6772.5
$ws.Range("M77").Value = -2404.5

$ws.Range("H86").Value = 50000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 50000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 50000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -52372

$ws.Range("H89").Value = 50000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 50000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 150000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -161856

$ws.Range("H116").Value = 1316.8334
$ws.Range("I116").Value = 1377.6
$ws.Range("J116").Value = 1013
$ws.Range("K116").Value = 1377.6
$ws.Range("L116").Value = 1013
$ws.Range("M116").Value = 916.4000000000001
$ws.Range("N116").Value = -5601

$ws.Range("H125").Value = 85300.44500000001
$ws.Range("J125").Value = 85300.44500000001
$ws.Range("L125").Value = 85300.44500000001
$ws.Range("N125").Value = -95140.44500000001

$ws.Range("H131").Value = 60249.668
$ws.Range("J131").Value = 60249.668
$ws.Range("L131").Value = 60249.668
$ws.Range("N131").Value = -70329.66800000001

$ws.Range("H132").Value = 2946.5217
$ws.Range("I132").Value = 2901.025
$ws.Range("K132").Value = 8703.075000000001
$ws.Range("M132").Value = -6173.075000000001

$ws.Range("H136").Value = 2401.7778
$ws.Range("I136").Value = 2396
$ws.Range("K136").Value = 7188
$ws.Range("M136").Value = -4638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1316.8334
$ws.Range("I3").Value = 1377.6
$ws.Range("J3").Value = 1013
$ws.Range("K3").Value = 1377.6
$ws.Range("L3").Value = 1013
$ws.Range("M3").Value = -1263.6
$ws.Range("N3").Value = -1241

$ws.Range("H7").Value = 1214.2858

$ws.Range("H20").Value = 3268
$ws.Range("I20").Value = 3048.2856
$ws.Range("J20").Value = 3687.4546
$ws.Range("K20").Value = 3048.2856
$ws.Range("L20").Value = 3687.4546
$ws.Range("M20").Value = -2801.2856
$ws.Range("N20").Value = -4181.4546

$ws.Range("H86").Value = 812121.4399999999
$ws.Range("I86").Value = 1419133.4
$ws.Range("K86").Value = 1419133.4
$ws.Range("M86").Value = -1418010.4

$ws.Range("H89").Value = 812121.4399999999
$ws.Range("I89").Value = 1419133.4
$ws.Range("K89").Value = 7095667
$ws.Range("M89").Value = -7090051

$ws.Range("H107").Value = 590765.25
$ws.Range("J107").Value = 2503264.2
$ws.Range("L107").Value = 2503264.2
$ws.Range("N107").Value = -2507104.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 147472.86
$ws.Range("I31").Value = 1828
$ws.Range("J31").Value = 341666
$ws.Range("K31").Value = 1828
$ws.Range("L31").Value = 341666
$ws.Range("M31").Value = -1533
$ws.Range("N31").Value = -342256

$ws.Range("H34").Value = 147472.86
$ws.Range("I34").Value = 1828
$ws.Range("J34").Value = 341666
$ws.Range("K34").Value = 1828
$ws.Range("L34").Value = 341666
$ws.Range("M34").Value = -1626
$ws.Range("N34").Value = -342070

$ws.Range("H35").Value = 8710
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H107").Value = 634.2857
$ws.Range("I107").Value = 516.26086
$ws.Range("K107").Value = 516.26086
$ws.Range("M107").Value = 1403.73914

$ws.Range("H140").Value = 42499
$ws.Range("J140").Value = 49998.8
$ws.Range("L140").Value = 49998.8
$ws.Range("N140").Value = -60358.8

$ws.Range("H141").Value = 300886.62
$ws.Range("J141").Value = 311509.6
$ws.Range("L141").Value = 311509.6
$ws.Range("N141").Value = -321869.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 19001.666
$ws.Range("I87").Value = 19001.666
$ws.Range("K87").Value = 57004.99800000001
$ws.Range("M87").Value = -55756.99800000001

$ws.Range("H90").Value = 19001.666
$ws.Range("I90").Value = 19001.666
$ws.Range("K90").Value = 171014.994
$ws.Range("M90").Value = -164774.994

$ws.Range("H114").Value = 765.4
$ws.Range("I114").Value = 374
$ws.Range("K114").Value = 1122
$ws.Range("M114").Value = 2132

$ws.Range("H118").Value = 21101.889
$ws.Range("I118").Value = 4009.6667
$ws.Range("J118").Value = 29648
$ws.Range("K118").Value = 12029.0001
$ws.Range("L118").Value = 88944
$ws.Range("M118").Value = -10786.0001
$ws.Range("N118").Value = -91430

$ws.Range("H131").Value = 2812.8
$ws.Range("J131").Value = 3128.9768
$ws.Range("L131").Value = 9386.930399999999
$ws.Range("N131").Value = -19466.9304

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7356.731
$ws.Range("I70").Value = 5987.647
$ws.Range("K70").Value = 5987.647
$ws.Range("M70").Value = -5717.647

$ws.Range("H73").Value = 7356.731
$ws.Range("I73").Value = 5987.647
$ws.Range("K73").Value = 5987.647
$ws.Range("M73").Value = -5051.647

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2983.35
$ws.Range("I61").Value = 3087.9473
$ws.Range("J61").Value = 996
$ws.Range("K61").Value = 3087.9473
$ws.Range("L61").Value = 996
$ws.Range("M61").Value = -2885.9473
$ws.Range("N61").Value = -1400

$ws.Range("H113").Value = 2983.35
$ws.Range("I113").Value = 3087.9473
$ws.Range("J113").Value = 996
$ws.Range("K113").Value = 3087.9473
$ws.Range("L113").Value = 996
$ws.Range("M113").Value = -917.9472999999998
$ws.Range("N113").Value = -5336

$ws.Range("H136").Value = 255778.72
$ws.Range("I136").Value = 458726.97
$ws.Range("K136").Value = 1376180.91
$ws.Range("M136").Value = -1373630.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1369.8462
$ws.Range("I107").Value = 1550.3334
$ws.Range("J107").Value = 611.8
$ws.Range("K107").Value = 4651.0002
$ws.Range("L107").Value = 1835.4
$ws.Range("M107").Value = -2731.0002
$ws.Range("N107").Value = -5675.4

$ws.Range("H126").Value = 1750.6666
$ws.Range("I126").Value = 1633
$ws.Range("J126").Value = 1868.3334
$ws.Range("K126").Value = 4899
$ws.Range("L126").Value = 5605.0002
$ws.Range("M126").Value = -2429
$ws.Range("N126").Value = -10545.0002
